$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C for the "navMesh" field (old C/resPath, D/type,
# E/position shift right to D/E/F).
$ws.Columns("C:C").Insert()

# --- Row 1: English field keys, now tagged with (C)/(S) export markers ---
$ws.Range("A1").Value = "mapID(I)"
$ws.Range("B1").Value = "(C)mapName(S)"
$ws.Range("C1").Value = "(C)navMesh(S)"
$ws.Range("D1").Value = "(S)resPath(S)"
$ws.Range("E1").Value = "type(I)"
$ws.Range("F1").Value = "position(V)"

# --- Row 2: Chinese display labels ---
$ws.Range("A2").Value = "地图ID"
$ws.Range("B2").Value = "地图名称(客户端)"
$ws.Range("C2").Value = "navMesh文件（客户端）"
$ws.Range("D2").Value = "地图信息（服务器）"
$ws.Range("E2").Value = "类型"
$ws.Range("F2").Value = "出生点"

# --- Row 3: sample data values ---
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "scene/newplayerscene1"
$ws.Range("C3").Value = "newPlayerScene1"
$ws.Range("D3").Value = "spaces/newplayerscene1"
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = "0,0,0"

# Move the reviewer comment that used to sit on the "type" column (D2)
# onto its new location (E2) now that the column shifted right.
$oldComment = $ws.Range("D2").Comment
if ($oldComment) {
    $oldComment.Delete()
}
$ws.Range("E2").AddComment("作者:`n0:场景`n1:副本`n")

# Re-fit the data columns to their new (generally longer) content. The
# position/spawn-point column (F, formerly E) keeps its old width since its
# own content never changed - only its column index shifted.
$ws.Columns("B:D").AutoFit()
$ws.Columns("B:B").ColumnWidth = 23.15
$ws.Columns("C:C").ColumnWidth = 22.15
$ws.Columns("D:D").ColumnWidth = 24.24

# Restore the active selection used in the saved workbook.
$ws.Range("B2").Select() | Out-Null
